# Powerpoint writer: consolidate text run nodes.
# Merge adjacent "word" + "space" runs into a single run, leaving the
# following run untouched, by rewriting the text of a Characters()
# sub-range that spans exactly the runs to be merged.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title 1: "A" + " " + "slide"  ->  "A " + "slide" ---
$titleTr = $s.Shapes.Item(1).TextFrame.TextRange
$titleTr.Characters(1, 2).Text = "A "

# --- Table (Content Placeholder 5), cell (1,2): "a" + " " + "table" -> "a " + "table" ---
$tbl = $s.Shapes.Item(3).Table
$cellTr = $tbl.Cell(1, 2).Shape.TextFrame.TextRange
$cellTr.Characters(1, 2).Text = "a "

# --- TextBox 3: "Plus" + " " + "an" + " " + "image" -> "Plus " + "an " + "image" ---
$tbTr = $s.Shapes.Item(7).TextFrame.TextRange
$tbTr.Characters(1, 5).Text = "Plus "
$tbTr.Characters(6, 3).Text = "an "
